$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 5829.55
$ws.Range("I18").Value = 5829.55
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 5829.55
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -5545.55
$ws.Range("N18").ClearContents()

# Row 40
$ws.Range("H40").Value = 5840.5454
$ws.Range("I40").Value = 4666.5
$ws.Range("J40").Value = 7249.4
$ws.Range("K40").Value = 4666.5
$ws.Range("L40").Value = 7249.4
$ws.Range("M40").Value = -4491.5
$ws.Range("N40").Value = -7599.4

# Row 46
$ws.Range("H46").Value = 115177.336
$ws.Range("J46").Value = 147670.86
$ws.Range("L46").Value = 443012.58
$ws.Range("N46").Value = -443250.58

# Row 54
$ws.Range("H54").Value = 11250
$ws.Range("J54").Value = 12000
$ws.Range("L54").Value = 12000
$ws.Range("N54").Value = -12972

# Row 59
$ws.Range("I59").Value = 1700
$ws.Range("J59").Value = 1246.6666
$ws.Range("K59").Value = 5100
$ws.Range("L59").Value = 3739.9998
$ws.Range("M59").Value = -4543
$ws.Range("N59").Value = -4853.9998

# Row 60
$ws.Range("H60").Value = 115177.336
$ws.Range("J60").Value = 147670.86
$ws.Range("L60").Value = 443012.58
$ws.Range("N60").Value = -443980.58

# Row 113
$ws.Range("H113").Value = 8411.076999999999
$ws.Range("I113").Value = 9536.75
$ws.Range("K113").Value = 9536.75
$ws.Range("M113").Value = -6282.75

# Row 132
$ws.Range("H132").Value = 2266.9688
$ws.Range("I132").Value = 2148.9673
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 6446.901899999999
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -3916.901899999999
$ws.Range("N132").Value = -19059.0005

# Row 133
$ws.Range("H133").Value = 99777
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 99777
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 99777
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -109897

# Row 136
$ws.Range("H136").Value = 73167
$ws.Range("J136").Value = 73167
$ws.Range("L136").Value = 73167
$ws.Range("N136").Value = -83367

# Row 138
$ws.Range("H138").Value = 3958.4922
$ws.Range("I138").Value = 1845.8182
$ws.Range("J138").Value = 5039.3955
$ws.Range("K138").Value = 5537.4546
$ws.Range("L138").Value = 15118.1865
$ws.Range("M138").Value = -397.4546
$ws.Range("N138").Value = -25398.1865

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6995.25
$ws.Range("I45").Value = 7030
$ws.Range("J45").Value = 6968.222
$ws.Range("K45").Value = 7030
$ws.Range("L45").Value = 6968.222
$ws.Range("M45").Value = -6653
$ws.Range("N45").Value = -7722.222

# Row 61
$ws.Range("H61").Value = 6457.4614
$ws.Range("I61").Value = 6663.75
$ws.Range("K61").Value = 6663.75
$ws.Range("M61").Value = -6451.75

# Row 74
$ws.Range("H74").Value = 6489.8335
$ws.Range("I74").Value = 1775.6
$ws.Range("K74").Value = 1775.6
$ws.Range("M74").Value = -901.5999999999999

# Row 77
$ws.Range("H77").Value = 6489.8335
$ws.Range("I77").Value = 1775.6
$ws.Range("K77").Value = 8878
$ws.Range("M77").Value = -4510

# Row 136
$ws.Range("H136").Value = 6457.4614
$ws.Range("I136").Value = 6663.75
$ws.Range("K136").Value = 19991.25
$ws.Range("M136").Value = -17441.25

$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 5750
$ws.Range("I29").Value = 1500
$ws.Range("K29").Value = 1500
$ws.Range("M29").Value = -1211

# Row 100
$ws.Range("H100").Value = 33908
$ws.Range("J100").Value = 33908
$ws.Range("L100").Value = 33908
$ws.Range("N100").Value = -36072

# Row 105
$ws.Range("H105").Value = 3159.2285
$ws.Range("I105").Value = 2507.6667
$ws.Range("J105").Value = 5358.25
$ws.Range("K105").Value = 2507.6667
$ws.Range("L105").Value = 5358.25
$ws.Range("M105").Value = -760.6667000000002
$ws.Range("N105").Value = -8852.25

# Row 134
$ws.Range("H134").Value = 16170.5
$ws.Range("I134").Value = 18700.215
$ws.Range("J134").Value = 7316.5
$ws.Range("K134").Value = 56100.645
$ws.Range("L134").Value = 21949.5
$ws.Range("M134").Value = -53565.645
$ws.Range("N134").Value = -27019.5

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 250

# Row 132
$ws.Range("H132").Value = 21748.4
$ws.Range("I132").Value = 1325.3889
$ws.Range("J132").Value = 205555.5
$ws.Range("K132").Value = 3976.1667
$ws.Range("L132").Value = 616666.5
$ws.Range("M132").Value = -1446.1667
$ws.Range("N132").Value = -621726.5

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 2761.2
$ws.Range("J113").Value = 2761.2
$ws.Range("L113").Value = 8283.599999999999
$ws.Range("N113").Value = -12623.6

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 4166.3335
$ws.Range("J21").Value = 4166.3335
$ws.Range("L21").Value = 4166.3335
$ws.Range("N21").Value = -4512.3335

# Row 30
$ws.Range("H30").Value = 4166.3335
$ws.Range("J30").Value = 4166.3335
$ws.Range("L30").Value = 4166.3335
$ws.Range("N30").Value = -4376.3335

# Row 102
$ws.Range("H102").Value = 24487.834
$ws.Range("I102").Value = 27885.4
$ws.Range("J102").Value = 7500
$ws.Range("K102").Value = 27885.4
$ws.Range("L102").Value = 7500
$ws.Range("M102").Value = -26263.4
$ws.Range("N102").Value = -10744

# Row 122
$ws.Range("H122").Value = 8838.893
$ws.Range("I122").Value = 6635.316
$ws.Range("K122").Value = 19905.948
$ws.Range("M122").Value = -17455.948

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 144285710
$ws.Range("J2").Value = 32000000
$ws.Range("L2").Value = 32000000
$ws.Range("N2").Value = -32000224

# Row 122
$ws.Range("H122").Value = 7629.206
$ws.Range("I122").Value = 5565.3447
$ws.Range("J122").Value = 19599.6
$ws.Range("K122").Value = 16696.0341
$ws.Range("L122").Value = 58798.8
$ws.Range("M122").Value = -14246.0341
$ws.Range("N122").Value = -63698.8

# Row 136
$ws.Range("H136").Value = 6883.64
$ws.Range("I136").Value = 2064.8235
$ws.Range("J136").Value = 17123.625
$ws.Range("K136").Value = 6194.470499999999
$ws.Range("L136").Value = 51370.875
$ws.Range("M136").Value = -3644.470499999999
$ws.Range("N136").Value = -56470.875

$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 1436572.9
$ws.Range("I20").Value = 6010
$ws.Range("J20").Value = 1675000
$ws.Range("K20").Value = 6010
$ws.Range("L20").Value = 1675000
$ws.Range("M20").Value = -5770
$ws.Range("N20").Value = -1675480

# Row 96
$ws.Range("H96").Value = 20002998
$ws.Range("I96").Value = 33335332
$ws.Range("J96").Value = 4500
$ws.Range("K96").Value = 33335332
$ws.Range("L96").Value = 4500
$ws.Range("M96").Value = -33333959
$ws.Range("N96").Value = -7246

# Row 122
$ws.Range("H122").Value = 5081.2085
$ws.Range("I122").Value = 3035.3076
$ws.Range("K122").Value = 9105.9228
$ws.Range("M122").Value = -6655.9228

# Row 126
$ws.Range("H126").Value = 35393.617
$ws.Range("I126").Value = 42512.7
$ws.Range("J126").Value = 11663.333
$ws.Range("K126").Value = 127538.1
$ws.Range("L126").Value = 34989.999
$ws.Range("M126").Value = -125068.1
$ws.Range("N126").Value = -39929.999

# Row 132
$ws.Range("H132").Value = 11630
$ws.Range("I132").Value = 11730.833
$ws.Range("K132").Value = 35192.499
$ws.Range("M132").Value = -32662.499

# Row 135
$ws.Range("H135").Value = 147500
$ws.Range("J135").Value = 147500
$ws.Range("L135").Value = 147500
$ws.Range("N135").Value = -157640
